# Matlab scripts for diagrams
# Adds the "miniboone pid" benchmark results (columns E-K) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("E1").Value = "miniboone pid"
$ws.Range("F1").Value = "bbc"
$ws.Range("G1").Value = "cnn"
$ws.Range("H1").Value = "cnnibn"
$ws.Range("I1").Value = "ndtv"
$ws.Range("J1").Value = "timesnow"
$ws.Range("K1").Value = "features"
$ws.Range("E1:K1").Style = "Bad"

# --- Data rows 2-6 (full set of new columns) -------------------------------
$ws.Range("E2").Value = 49515.9
$ws.Range("F2").Value = 332.744
$ws.Range("G2").Value = 576.034
$ws.Range("H2").Value = 1038.02
$ws.Range("I2").Value = 809.125
$ws.Range("J2").Value = 2212.08
$ws.Range("K2").Value = 484865

$ws.Range("E3").Value = 31458.5
$ws.Range("F3").Value = 246.289
$ws.Range("G3").Value = 376.915
$ws.Range("H3").Value = 614.14
$ws.Range("I3").Value = 491.827
$ws.Range("J3").Value = 1369.55
$ws.Range("K3").Value = 310285

$ws.Range("E4").Value = 24016.13
$ws.Range("F4").Value = 181.389
$ws.Range("G4").Value = 265.583
$ws.Range("H4").Value = 432.348
$ws.Range("I4").Value = 313.382
$ws.Range("J4").Value = 903.358
$ws.Range("K4").Value = 206237

$ws.Range("E5").Value = 20169.15
$ws.Range("F5").Value = 147.338
$ws.Range("G5").Value = 213.567
$ws.Range("H5").Value = 319.798
$ws.Range("I5").Value = 243.445
$ws.Range("J5").Value = 650.778
$ws.Range("K5").Value = 133446

$ws.Range("E6").Value = 19419.1
$ws.Range("F6").Value = 101.142
$ws.Range("G6").Value = 159.161
$ws.Range("H6").Value = 247.616
$ws.Range("I6").Value = 182.167
$ws.Range("J6").Value = 480.224
$ws.Range("K6").Value = 91637

# --- Data rows 7-9 (only column E has data) --------------------------------
$ws.Range("E7").Value = 13105.2
$ws.Range("E8").Value = 11941.7
$ws.Range("E9").Value = 9151.75

# --- Column E width ----------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 13.5

# --- View / selection state ---------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 160
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("K7").Select() | Out-Null
$win.Left = -110
$win.Top = -110
$win.Width = 19420
$win.Height = 11020
